{"js": "// Replace the division-problem text in the document's table cells.\n// Each \"before\" string is unique in the document, so a scoped\n// search-and-replace (matchCase, on the whole body) is safe and will\n// only ever touch the intended cell's run.\nconst replacements = [\n  [\"59\u00f78=\", \"27\u00f75=\"],\n  [\"45\u00f77=\", \"42\u00f73=\"],\n  [\"34\u00f78=\", \"30\u00f74=\"],\n  [\"38\u00f72=\", \"38\u00f75=\"],\n  [\"50\u00f73=\", \"88\u00f75=\"],\n  [\"45\u00f75=\", \"23\u00f73=\"],\n  [\"47\u00f72=\", \"38\u00f79=\"],\n  [\"56\u00f73=\", \"37\u00f78=\"],\n  [\"10\u00f77=\", \"96\u00f79=\"],\n  [\"32\u00f76=\", \"46\u00f78=\"],\n  [\"26\u00f76=\", \"27\u00f72=\"],\n  [\"29\u00f75=\", \"24\u00f76=\"],\n  [\"44\u00f79=\", \"45\u00f74=\"],\n  [\"66\u00f72=\", \"37\u00f75=\"],\n  [\"14\u00f72=\", \"86\u00f72=\"],\n  [\"35\u00f76=\", \"18\u00f72=\"],\n  [\"79\u00f77=\", \"16\u00f77=\"],\n  [\"95\u00f74=\", \"83\u00f76=\"],\n  [\"77\u00f77=\", \"47\u00f74=\"],\n  [\"38\u00f77=\", \"12\u00f73=\"],\n  [\"41\u00f78=\", \"20\u00f74=\"],\n  [\"90\u00f78=\", \"82\u00f76=\"],\n  [\"90\u00f74=\", \"38\u00f79=\"],\n  [\"43\u00f77=\", \"91\u00f78=\"],\n  [\"23\u00f72=\", \"26\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in the document's table cells.\n# Each \"before\" string is unique in the document, so a document-wide\n# Find/Replace (MatchCase) is safe and only touches the intended cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"59\u00f78=\", \"27\u00f75=\"),\n    @(\"45\u00f77=\", \"42\u00f73=\"),\n    @(\"34\u00f78=\", \"30\u00f74=\"),\n    @(\"38\u00f72=\", \"38\u00f75=\"),\n    @(\"50\u00f73=\", \"88\u00f75=\"),\n    @(\"45\u00f75=\", \"23\u00f73=\"),\n    @(\"47\u00f72=\", \"38\u00f79=\"),\n    @(\"56\u00f73=\", \"37\u00f78=\"),\n    @(\"10\u00f77=\", \"96\u00f79=\"),\n    @(\"32\u00f76=\", \"46\u00f78=\"),\n    @(\"26\u00f76=\", \"27\u00f72=\"),\n    @(\"29\u00f75=\", \"24\u00f76=\"),\n    @(\"44\u00f79=\", \"45\u00f74=\"),\n    @(\"66\u00f72=\", \"37\u00f75=\"),\n    @(\"14\u00f72=\", \"86\u00f72=\"),\n    @(\"35\u00f76=\", \"18\u00f72=\"),\n    @(\"79\u00f77=\", \"16\u00f77=\"),\n    @(\"95\u00f74=\", \"83\u00f76=\"),\n    @(\"77\u00f77=\", \"47\u00f74=\"),\n    @(\"38\u00f77=\", \"12\u00f73=\"),\n    @(\"41\u00f78=\", \"20\u00f74=\"),\n    @(\"90\u00f78=\", \"82\u00f76=\"),\n    @(\"90\u00f74=\", \"38\u00f79=\"),\n    @(\"43\u00f77=\", \"91\u00f78=\"),\n    @(\"23\u00f72=\", \"26\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $after\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($before, $true, $false, $false, $false, $false, $true, 1, $false, $after, 2)\n}\n"}
